# UndoRedoStack: update to store ReadOnlyAddressBook
#
# Sequence-diagram changes on the slide:
#  1. The participant box "x:XYZCommand" becomes ":ModelManager".
#  2. The message "undo()" (sent to the stack) becomes
#     "resetData(AddressBook)", and its textbox is widened/repositioned
#     to fit the new label.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Rectangle 62 / "x:XYZCommand" -> ":ModelManager" -------------------
$xyz = $s.Shapes.Item(23)
$xyzTr = $xyz.TextFrame.TextRange
$xyzTr.Text = ":ModelManager"
# Touch the two logical runs (":" and "ModelManager") individually so they
# stay separate runs in the saved XML, the way PowerPoint keeps a
# mis-spelling-flagged word in its own run.
$xyzTr.Characters(1, 1).Font.Size = 16
$xyzTr.Characters(2, 12).Font.Size = 16

# --- 2. TextBox 87 / "undo()" -> "resetData(AddressBook)" ------------------
$undoShp = $s.Shapes.Item(38)

# Reposition/resize the textbox to its new extents.
# (EMU -> points, 12700 EMU per point; target EMU: off 6142472,3416411 ext 1615428,184666)
$undoShp.Left = 483.6592325984252
$undoShp.Top = 269.0087501574803
$undoShp.Width = 127.19905511811024
$undoShp.Height = 14.540629921259843

$undoTr = $undoShp.TextFrame.TextRange
$undoTr.Text = "resetData(AddressBook)"
# "resetData" (1-9) + "(" (10) + "AddressBook" (11-21) + ")" (22)
$undoTr.Characters(1, 9).Font.Size = 12
$undoTr.Characters(10, 1).Font.Size = 12
$undoTr.Characters(11, 11).Font.Size = 12
$undoTr.Characters(22, 1).Font.Size = 12
